# Apply the updated test-results numbers (re-run of evaluation with new thresholds)
# to the "miscellaneous folder" workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Second table (I:O, rows 2-9) -- Initial/Final/Max IoU updates ----
$ws.Range("J2").Value = 0.63
$ws.Range("K2").Value = 0.53
$ws.Range("L2").Value = 0.63

$ws.Range("J3").Value = 0.34
$ws.Range("K3").Value = 0.21
$ws.Range("L3").Value = 0.34
$ws.Range("M3").Value = 10
$ws.Range("O3").Value = "max number of steps"

$ws.Range("J4").Value = 0.15
$ws.Range("K4").Value = 0.11
$ws.Range("L4").Value = 0.15

# I5 becomes a numeric-looking image id ("007857") -- go through a formula +
# paste-values round trip so Excel keeps it as text instead of silently
# coercing it to the number 7857 (no leading zeros, no <f>).
$ws.Range("I5").Formula = '="007857"'
$ws.Range("I5").Copy()
$ws.Range("I5").PasteSpecial(-4163)   # xlPasteValues
$ws.Application.CutCopyMode = $false

$ws.Range("J5").Value = 0.47
$ws.Range("K5").Value = 0.45
$ws.Range("L5").Value = 0.47
$ws.Range("M5").Value = 4
$ws.Range("O5").Value = "STOP q-value"

$ws.Range("J6").Value = 0.86
$ws.Range("K6").Value = 0.83
$ws.Range("L6").Value = 0.87

$ws.Range("J7").Value = 0.91
$ws.Range("K7").Value = 0.78
$ws.Range("L7").Value = 0.91
$ws.Range("M7").Value = 3

$ws.Range("J8").Value = 0.71
$ws.Range("K8").Value = 0.55000000000000004
$ws.Range("L8").Value = 0.71
$ws.Range("M8").Value = 7
$ws.Range("O8").Value = "STOP q-value"

$ws.Range("J9").Value = 0.46
$ws.Range("K9").Value = 0.43
$ws.Range("L9").Value = 0.46
$ws.Range("M9").Value = 3

# ---- First table (A:G, rows 2-19) ----
# Row 5 and row 6 image-name swap ("007857" <-> "007820")
$ws.Range("A6").Formula = '="007820"'
$ws.Range("A6").Copy()
$ws.Range("A6").PasteSpecial(-4163)   # xlPasteValues
$ws.Application.CutCopyMode = $false

$ws.Range("B6").Value = 0.33
$ws.Range("C6").Value = 0.53
$ws.Range("D6").Value = 0.53
$ws.Range("E6").Value = 10
$ws.Range("G6").Value = "max number of steps"

$ws.Range("B14").Value = 0.41

# ---- Summary counters ----
$ws.Range("B23").Value = 4
$ws.Range("J23").Value = 4

# ---- Footer block: drop row 29, shift the 3 footer lines down one row ----
$ws.Range("C29").ClearContents()
$ws.Range("C30").Value = "IoU increase: 18"
$ws.Range("C31").Value = "IoU decrease: 8"
$ws.Range("C32").Value = "IoU average: 0.5635587432643805"

# ---- View state: scroll + selection ----
$ws.Range("C30:C32").Select()
